# "Generate Report for Archive"
#
# The localization status report is regenerated: the "Status" value that
# used to read "Ready for handoff" is now "In Translation" (it shows up on
# the Overview sheet for both locale columns, and on each locale sheet's
# own Status column). The Status column also got narrower on every sheet
# that has one.

$wb = $excel.ActiveWorkbook

# --- 1. Update the Status text wherever it appears -------------------------
$overview = $wb.Sheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Sheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Sheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# --- 2. Narrow the Status columns ------------------------------------------
# Raw OOXML column width goes from 17.2159881591797 down to 13.4101845877511
# (character units, i.e. ColumnWidth + ~0.8333 offset applied by the host).
$newWidth = 13.4101845877511 - 0.8333333333333334

$overview.Columns.Item(5).ColumnWidth = $newWidth
$overview.Columns.Item(6).ColumnWidth = $newWidth
$zhcn.Columns.Item(3).ColumnWidth = $newWidth
$dede.Columns.Item(3).ColumnWidth = $newWidth
